# case 1 data update: row 1 grows from 11 cols (A:K) to 17 cols (A:Q)
# and every existing value is replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths -------------------------------------------------
# Column F narrows from the "3.14" group down to the "2.14" group;
# columns I:L are new members of the "3.14" group; columns M:Q are new
# members of the "5.71" group (same widths already used by I:K/A:B).
# (ColumnWidth is quantised to the nearest pixel by the host, so these
# inputs are chosen to land closest to the 2.140625 / 3.140625 /
# 5.7109375 character-width targets.)
$ws.Columns.Item(6).ColumnWidth = 1.3

$ws.Columns.Item(9).ColumnWidth = 2.3
$ws.Columns.Item(10).ColumnWidth = 2.3
$ws.Columns.Item(11).ColumnWidth = 2.3
$ws.Columns.Item(12).ColumnWidth = 2.3

$ws.Columns.Item(13).ColumnWidth = 4.8
$ws.Columns.Item(14).ColumnWidth = 4.8
$ws.Columns.Item(15).ColumnWidth = 4.8
$ws.Columns.Item(16).ColumnWidth = 4.8
$ws.Columns.Item(17).ColumnWidth = 4.8

# --- row 1 values ----------------------------------------------------
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 32
$ws.Range("D1").Value = 23
$ws.Range("E1").Value = 16
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 28
$ws.Range("H1").Value = 16
$ws.Range("I1").Value = 12
$ws.Range("J1").Value = 20
$ws.Range("K1").Value = 32
$ws.Range("L1").Value = 31
$ws.Range("M1").Value = 0.090999999999999998
$ws.Range("N1").Value = 0.048999999999999995
$ws.Range("O1").Value = 0.014000000000000002
$ws.Range("P1").Value = 0.063
$ws.Range("Q1").Value = 0.023999999999999997
